$wb = $excel.ActiveWorkbook

# --- Sheet1: Neg_Change ---
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws1.Range("A2").Value = "ULTRACEMCO"
$ws1.Range("B2").Value = 12040
$ws1.Range("C2").Value = 12084
$ws1.Range("D2").Value = 11950
$ws1.Range("E2").Value = 12040
$ws1.Range("F2").Value = 159943
$ws1.Range("G2").Value = 394272
$ws1.Range("H2").Value = -0.5943333536238942
$ws1.Range("I2").Value = "ULTRACEMCO"

$ws1.Range("A3").Value = "TATASTEEL"
$ws1.Range("B3").Value = 184.4
$ws1.Range("C3").Value = 185.19
$ws1.Range("D3").Value = 182.64
$ws1.Range("E3").Value = 184.45
$ws1.Range("F3").Value = 24270977
$ws1.Range("G3").Value = 51610530
$ws1.Range("H3").Value = -0.5297281969396556
$ws1.Range("I3").Value = "TATASTEEL"

$ws1.Range("A4").Value = "NTPC"
$ws1.Range("B4").Value = 349
$ws1.Range("C4").Value = 349
$ws1.Range("D4").Value = 343.45
$ws1.Range("E4").Value = 345.3
$ws1.Range("F4").Value = 6721032
$ws1.Range("G4").Value = 15423201
$ws1.Range("H4").Value = -0.5642258698437503
$ws1.Range("I4").Value = "NTPC"

$ws1.Range("A5").Value = "BOSCHLTD"
$ws1.Range("B5").Value = 37500
$ws1.Range("C5").Value = 37570
$ws1.Range("D5").Value = 36970
$ws1.Range("E5").Value = 37240
$ws1.Range("F5").Value = 24286
$ws1.Range("G5").Value = 49599
$ws1.Range("H5").Value = -0.5103530313111151
$ws1.Range("I5").Value = "BOSCHLTD"

$ws1.Range("A6").Value = "ICICIGI"
$ws1.Range("B6").Value = 2029.9
$ws1.Range("C6").Value = 2032.8
$ws1.Range("D6").Value = 2005
$ws1.Range("E6").Value = 2015
$ws1.Range("F6").Value = 213019
$ws1.Range("G6").Value = 420554
$ws1.Range("H6").Value = -0.4934800287240164
$ws1.Range("I6").Value = "ICICIGI"

$ws1.Range("A7").Value = "LICI"
$ws1.Range("B7").Value = 910.85
$ws1.Range("C7").Value = 910.85
$ws1.Range("D7").Value = 897.1
$ws1.Range("E7").Value = 900.5
$ws1.Range("F7").Value = 756385
$ws1.Range("G7").Value = 1649697
$ws1.Range("H7").Value = -0.5415006513317294
$ws1.Range("I7").Value = "LICI"

$ws1.Range("A8").Value = "NHPC"
$ws1.Range("B8").Value = 87.14
$ws1.Range("C8").Value = 87.23
$ws1.Range("D8").Value = 85.70999999999999
$ws1.Range("E8").Value = 86.2
$ws1.Range("F8").Value = 7400279
$ws1.Range("G8").Value = 17056980
$ws1.Range("H8").Value = -0.5661436549729202
$ws1.Range("I8").Value = "NHPC"

$ws1.Range("A9").Value = "CAMS"
$ws1.Range("B9").Value = 3874.3
$ws1.Range("C9").Value = 3980
$ws1.Range("D9").Value = 3860
$ws1.Range("E9").Value = 3951
$ws1.Range("F9").Value = 663840
$ws1.Range("G9").Value = 1373940
$ws1.Range("H9").Value = -0.5168347962793135
$ws1.Range("I9").Value = "CAMS"

$ws1.Range("A10").Value = "ANGELONE"
$ws1.Range("B10").Value = 2513
$ws1.Range("C10").Value = 2532.9
$ws1.Range("D10").Value = 2495.1
$ws1.Range("E10").Value = 2513
$ws1.Range("F10").Value = 324753
$ws1.Range("G10").Value = 716950
$ws1.Range("H10").Value = -0.5470353581142339
$ws1.Range("I10").Value = "ANGELONE"

$ws1.Range("A11").Value = "CDSL"
$ws1.Range("B11").Value = 1613.8
$ws1.Range("C11").Value = 1625
$ws1.Range("D11").Value = 1607
$ws1.Range("E11").Value = 1613.7
$ws1.Range("F11").Value = 809134
$ws1.Range("G11").Value = 1844623
$ws1.Range("H11").Value = -0.561355355538774
$ws1.Range("I11").Value = "CDSL"

$ws1.Range("A12").Value = "PGEL"
$ws1.Range("B12").Value = 573.2
$ws1.Range("C12").Value = 579.85
$ws1.Range("D12").Value = 571.5
$ws1.Range("E12").Value = 573.15
$ws1.Range("F12").Value = 575389
$ws1.Range("G12").Value = 1259123
$ws1.Range("H12").Value = -0.5430239936844931
$ws1.Range("I12").Value = "PGEL"

$ws1.Range("A13").Value = "MCX"
$ws1.Range("B13").Value = 9150.5
$ws1.Range("C13").Value = 9173.5
$ws1.Range("D13").Value = 9061
$ws1.Range("E13").Value = 9096.5
$ws1.Range("F13").Value = 194394
$ws1.Range("G13").Value = 419900
$ws1.Range("H13").Value = -0.5370469159323649
$ws1.Range("I13").Value = "MCX"

$ws1.Range("A14").Value = "POONAWALLA"
$ws1.Range("B14").Value = 487.8
$ws1.Range("C14").Value = 488.35
$ws1.Range("D14").Value = 480.4
$ws1.Range("E14").Value = 482.75
$ws1.Range("F14").Value = 681219
$ws1.Range("G14").Value = 1654852
$ws1.Range("H14").Value = -0.5883504990174347
$ws1.Range("I14").Value = "POONAWALLA"

# Remove rows 15-19 (no longer present in target data)
$ws1.Range("A15:I19").Delete()

# --- Sheet2: Pos_Change ---
$ws2 = $wb.Worksheets.Item("Pos_Change")
$ws2.Range("A2").Value = "BAJAJFINSV"
$ws2.Range("B2").Value = 2140
$ws2.Range("C2").Value = 2146
$ws2.Range("D2").Value = 2101.3
$ws2.Range("E2").Value = 2118
$ws2.Range("F2").Value = 488305
$ws2.Range("G2").Value = 343434
$ws2.Range("H2").Value = 0.4218306865365689
$ws2.Range("I2").Value = "BAJAJFINSV"

$ws2.Range("A3").Value = "BHARTIARTL"
$ws2.Range("B3").Value = 2070
$ws2.Range("C3").Value = 2088
$ws2.Range("D3").Value = 2058.6
$ws2.Range("E3").Value = 2069
$ws2.Range("F3").Value = 4222189
$ws2.Range("G3").Value = 2742750
$ws2.Range("H3").Value = 0.5393998723908486
$ws2.Range("I3").Value = "BHARTIARTL"

$ws2.Range("A4").Value = "BPCL"
$ws2.Range("B4").Value = 350
$ws2.Range("C4").Value = 359.8
$ws2.Range("D4").Value = 345.3
$ws2.Range("E4").Value = 358.1
$ws2.Range("F4").Value = 14807100
$ws2.Range("G4").Value = 9756868
$ws2.Range("H4").Value = 0.5176079045037814
$ws2.Range("I4").Value = "BPCL"

$ws2.Range("A5").Value = "TATAELXSI"
$ws2.Range("B5").Value = 5555
$ws2.Range("C5").Value = 5587
$ws2.Range("D5").Value = 5502.5
$ws2.Range("E5").Value = 5540
$ws2.Range("F5").Value = 123485
$ws2.Range("G5").Value = 81128
$ws2.Range("H5").Value = 0.5221008776254807
$ws2.Range("I5").Value = "TATAELXSI"

$ws2.Range("A6").Value = "FEDERALBNK"
$ws2.Range("B6").Value = 235.01
$ws2.Range("C6").Value = 237.44
$ws2.Range("D6").Value = 234.4
$ws2.Range("E6").Value = 234.65
$ws2.Range("F6").Value = 11944357
$ws2.Range("G6").Value = 7589866
$ws2.Range("H6").Value = 0.5737243582429519
$ws2.Range("I6").Value = "FEDERALBNK"

$ws2.Range("A7").Value = "VOLTAS"
$ws2.Range("B7").Value = 1410
$ws2.Range("C7").Value = 1425.2
$ws2.Range("D7").Value = 1400
$ws2.Range("E7").Value = 1415.3
$ws2.Range("F7").Value = 941460
$ws2.Range("G7").Value = 613025
$ws2.Range("H7").Value = 0.5357611842910158
$ws2.Range("I7").Value = "VOLTAS"

$ws2.Range("A8").Value = "EXIDEIND"
$ws2.Range("B8").Value = 385.3
$ws2.Range("C8").Value = 385.85
$ws2.Range("D8").Value = 377.6
$ws2.Range("E8").Value = 383.2
$ws2.Range("F8").Value = 3586392
$ws2.Range("G8").Value = 2270512
$ws2.Range("H8").Value = 0.57955210102391
$ws2.Range("I8").Value = "EXIDEIND"

$ws2.Range("A9").Value = "YESBANK"
$ws2.Range("B9").Value = 22.75
$ws2.Range("C9").Value = 22.75
$ws2.Range("D9").Value = 22.2
$ws2.Range("E9").Value = 22.24
$ws2.Range("F9").Value = 80860184
$ws2.Range("G9").Value = 56148847
$ws2.Range("H9").Value = 0.4401040861978875
$ws2.Range("I9").Value = "YESBANK"

$ws2.Range("A10").Value = "SBICARD"
$ws2.Range("B10").Value = 912
$ws2.Range("C10").Value = 915
$ws2.Range("D10").Value = 883.1
$ws2.Range("E10").Value = 886
$ws2.Range("F10").Value = 1323903
$ws2.Range("G10").Value = 900250
$ws2.Range("H10").Value = 0.47059483476812
$ws2.Range("I10").Value = "SBICARD"

Write-Host "Edit complete"